$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C13").Value = "sup_bati_km2"
$ws.Range("E13").Value = "area of neihborhood covered by buildings in square kilometers"

$ws.Range("E13").Select()
